$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain decimal-looking price strings (e.g. "304.51") that
# Excel would otherwise auto-coerce to numbers (dropping trailing/leading
# zeros, switching to scientific notation, etc). Briefly force the cell to
# Text format while assigning, then restore the default "Normal" style so
# the cell lands back at its original (unstyled / General) appearance.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '24.379.19'
$ws.Range('E2').Value = '  -5.75%  '
# Row 3
Set-TextValue $ws.Range('D3') '1.633.24'
$ws.Range('E3').Value = '  -7.33%  '
# Row 4
$ws.Range('E4').Value = '  +0.22%  '
# Row 5
Set-TextValue $ws.Range('D5') '1.002'
$ws.Range('E5').Value = '  +0.28%  '
# Row 6
Set-TextValue $ws.Range('D6') '304.51'
# Row 7
Set-TextValue $ws.Range('D7') '0.3618'
# Row 9
Set-TextValue $ws.Range('D9') '0.3227'
$ws.Range('E9').Value = '  -10.55%  '
# Row 10
Set-TextValue $ws.Range('D10') '1.101'
# Row 11
Set-TextValue $ws.Range('D11') '0.06878'
$ws.Range('E11').Value = '  -10.69%  '
# Row 12
Set-TextValue $ws.Range('D12') '1.003'
$ws.Range('E12').Value = '  +0.36%  '
# Row 13
Set-TextValue $ws.Range('D13') '5.894'
$ws.Range('E13').Value = '  -8.99%  '
# Row 14
Set-TextValue $ws.Range('D14') '19.08'
$ws.Range('E14').Value = '  -12.12%  '
# Row 15
Set-TextValue $ws.Range('D15') '1.640.33'
$ws.Range('E15').Value = '  -6.98%  '
# Row 16
Set-TextValue $ws.Range('D16') '6.510'
$ws.Range('E16').Value = '  -8.15%  '
# Row 17
Set-TextValue $ws.Range('D17') '0.00001042'
# Row 18
Set-TextValue $ws.Range('D18') '0.06510'
$ws.Range('E18').Value = '  -4.22%  '
# Row 19
$ws.Range('E19').Value = '  +0.25%  '
# Row 20
Set-TextValue $ws.Range('D20') '76.34'
$ws.Range('E20').Value = '  -12.29%  '
# Row 21
Set-TextValue $ws.Range('D21') '15.69'
$ws.Range('E21').Value = '  -11.28%  '
# Row 22
Set-TextValue $ws.Range('D22') '5.866'
$ws.Range('E22').Value = '  -9.86%  '
# Row 23
Set-TextValue $ws.Range('D23') '11.93'
$ws.Range('E23').Value = '  -8.28%  '
# Row 24
Set-TextValue $ws.Range('D24') '24.342.83'
$ws.Range('E24').Value = '  -5.68%  '
# Row 25
Set-TextValue $ws.Range('D25') '2.396'
$ws.Range('E25').Value = '  -1.70%  '
# Row 26
Set-TextValue $ws.Range('D26') '2.316'
$ws.Range('E26').Value = '  -20.29%  '
# Row 27
Set-TextValue $ws.Range('D27') '143.85'
$ws.Range('E27').Value = '  -7.78%  '
# Row 28
Set-TextValue $ws.Range('D28') '18.55'
$ws.Range('E28').Value = '  -10.79%  '
# Row 29
Set-TextValue $ws.Range('D29') '1.817.25'
$ws.Range('E29').Value = '  -7.26%  '
# Row 30
Set-TextValue $ws.Range('D30') '123.88'
$ws.Range('E30').Value = '  -7.66%  '
# Row 31
Set-TextValue $ws.Range('D31') '1.101'
$ws.Range('E31').Value = '  -8.88%  '
# Row 32
Set-TextValue $ws.Range('D32') '4.075'
$ws.Range('E32').Value = '  -3.29%  '
# Row 33
Set-TextValue $ws.Range('D33') '5.585'
$ws.Range('E33').Value = '  -22.52%  '
# Row 34
Set-TextValue $ws.Range('D34') '0.08351'
$ws.Range('E34').Value = '  -4.82%  '
# Row 35
Set-TextValue $ws.Range('D35') '1.666'
$ws.Range('E35').Value = '  -7.77%  '
# Row 36
Set-TextValue $ws.Range('D36') '12.23'
$ws.Range('E36').Value = '  -14.54%  '
# Row 37
Set-TextValue $ws.Range('D37') '5.081'
$ws.Range('E37').Value = '  -11.46%  '
# Row 38
Set-TextValue $ws.Range('D38') '0.05981'
$ws.Range('E38').Value = '  -11.54%  '
# Row 39
Set-TextValue $ws.Range('D39') '0.02208'
$ws.Range('E39').Value = '  -11.76%  '
# Row 40
Set-TextValue $ws.Range('D40') '1.201'
$ws.Range('E40').Value = '  -7.37%  '
# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D41') '8.114'
$ws.Range('E41').Value = '  -13.48%  '
# Row 42
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D42') '0.2026'
$ws.Range('E42').Value = '  -10.86%  '
# Row 43
Set-TextValue $ws.Range('D43') '1.001'
$ws.Range('E43').Value = '  +0.26%  '
# Row 44
Set-TextValue $ws.Range('D44') '0.5815'
$ws.Range('E44').Value = '  -11.84%  '
# Row 45
Set-TextValue $ws.Range('D45') '3.704'
$ws.Range('E45').Value = '  -5.24%  '
# Row 46
Set-TextValue $ws.Range('D46') '12.42'
$ws.Range('E46').Value = '  -13.74%  '
# Row 47
Set-TextValue $ws.Range('D47') '0.5507'
$ws.Range('E47').Value = '  -13.44%  '
# Row 48
Set-TextValue $ws.Range('D48') '121.47'
$ws.Range('E48').Value = '  -8.28%  '
# Row 49
Set-TextValue $ws.Range('D49') '1.910'
$ws.Range('E49').Value = '  -12.19%  '
# Row 50
Set-TextValue $ws.Range('D50') '0.06869'
$ws.Range('E50').Value = '  -8.60%  '
# Row 51
Set-TextValue $ws.Range('D51') '73.00'
$ws.Range('E51').Value = '  -9.92%  '
